$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.850.10"
$ws.Range("E2").Value = "  +2.62%  "
$ws.Range("D3").Value = "1.663.38"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.67"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.514"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.44"
$ws.Range("E8").Value = "  +2.92%  "
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0620"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0878"
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("D12").Value = "1.895.82"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "1.664.27"
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.13"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.550"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.90"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "248.40"
$ws.Range("E17").Value = "  +5.75%  "
$ws.Range("D18").Value = "27.805.29"
$ws.Range("E18").Value = "  +2.61%  "
$ws.Range("D19").Value = "0.0₃0732"
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("E20").Value = "  -4.22%  "
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.46"
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.36"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.58"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.18"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("E30").Value = "  +5.78%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.33"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.13"
$ws.Range("E33").Value = "  -3.20%  "
$ws.Range("D34").Value = "1.413.99"
$ws.Range("E34").Value = "  -8.31%  "
$ws.Range("E35").Value = "  -5.52%  "
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.928"
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.579"
$ws.Range("E38").Value = "  -4.24%  "
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("E40").Value = "  -3.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.12"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("E43").Value = "  -6.29%  "
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.791"
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.806.00"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.23"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  -4.13%  "
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("E51").Value = "  -0.42%  "
